$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2163.9285
$ws.Range("I62").Value = 2435.9092
$ws.Range("J62").Value = 1166.6666
$ws.Range("K62").Value = 2435.9092
$ws.Range("L62").Value = 1166.6666
$ws.Range("M62").Value = -1811.9092
$ws.Range("N62").Value = -2414.6666
$ws.Range("H65").Value = 2163.9285
$ws.Range("I65").Value = 2435.9092
$ws.Range("J65").Value = 1166.6666
$ws.Range("K65").Value = 12179.546
$ws.Range("L65").Value = 5833.333000000001
$ws.Range("M65").Value = -9059.546
$ws.Range("N65").Value = -12073.333
$ws.Range("H69").Value = 2853.25
$ws.Range("I69").Value = 1913
$ws.Range("J69").Value = 3166.6667
$ws.Range("K69").Value = 5739
$ws.Range("L69").Value = 9500.000100000001
$ws.Range("M69").Value = -4865
$ws.Range("N69").Value = -11248.0001
$ws.Range("H72").Value = 2853.25
$ws.Range("I72").Value = 1913
$ws.Range("J72").Value = 3166.6667
$ws.Range("K72").Value = 17217
$ws.Range("L72").Value = 28500.0003
$ws.Range("M72").Value = -12849
$ws.Range("N72").Value = -37236.0003
$ws.Range("H76").Value = 5184.533
$ws.Range("I76").Value = 6011.4443
$ws.Range("K76").Value = 6011.4443
$ws.Range("M76").Value = -5696.4443
$ws.Range("H79").Value = 5184.533
$ws.Range("I79").Value = 6011.4443
$ws.Range("K79").Value = 6011.4443
$ws.Range("M79").Value = -4919.4443
$ws.Range("H86").Value = 2838.8
$ws.Range("I86").Value = 2841.8572
$ws.Range("J86").Value = 2831.6667
$ws.Range("K86").Value = 2841.8572
$ws.Range("L86").Value = 2831.6667
$ws.Range("M86").Value = -1718.8572
$ws.Range("N86").Value = -5077.6667
$ws.Range("H89").Value = 2838.8
$ws.Range("I89").Value = 2841.8572
$ws.Range("J89").Value = 2831.6667
$ws.Range("K89").Value = 14209.286
$ws.Range("L89").Value = 14158.3335
$ws.Range("M89").Value = -8593.286
$ws.Range("N89").Value = -25390.3335
$ws.Range("H137").Value = 1626.6957
$ws.Range("I137").Value = 1389.0588
$ws.Range("J137").Value = 2300
$ws.Range("K137").Value = 4167.1764
$ws.Range("L137").Value = 6900
$ws.Range("M137").Value = -1617.1764
$ws.Range("N137").Value = -12000
$ws.Range("H138").Value = 1758.6809
$ws.Range("I138").Value = 1313.0646
$ws.Range("J138").Value = 2622.0625
$ws.Range("K138").Value = 3939.1938
$ws.Range("L138").Value = 7866.1875
$ws.Range("M138").Value = 1200.8062
$ws.Range("N138").Value = -18146.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1611.8113
$ws.Range("I61").Value = 1496.5745
$ws.Range("J61").Value = 2514.5
$ws.Range("K61").Value = 1496.5745
$ws.Range("L61").Value = 2514.5
$ws.Range("M61").Value = -1284.5745
$ws.Range("N61").Value = -2938.5
$ws.Range("H92").Value = 75919.89999999999
$ws.Range("J92").Value = 75919.89999999999
$ws.Range("L92").Value = 75919.89999999999
$ws.Range("N92").Value = -80911.89999999999
$ws.Range("H125").Value = 40476.668
$ws.Range("J125").Value = 40476.668
$ws.Range("L125").Value = 40476.668
$ws.Range("N125").Value = -50316.668
$ws.Range("H136").Value = 1611.8113
$ws.Range("I136").Value = 1496.5745
$ws.Range("J136").Value = 2514.5
$ws.Range("K136").Value = 4489.7235
$ws.Range("L136").Value = 7543.5
$ws.Range("M136").Value = -1939.7235
$ws.Range("N136").Value = -12643.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 91025.125
$ws.Range("J92").Value = 91025.125
$ws.Range("L92").Value = 91025.125
$ws.Range("N92").Value = -96017.125
$ws.Range("H134").Value = 3826.7
$ws.Range("I134").Value = 4258.7144
$ws.Range("J134").Value = 3594.077
$ws.Range("K134").Value = 12776.1432
$ws.Range("L134").Value = 10782.231
$ws.Range("M134").Value = -10241.1432
$ws.Range("N134").Value = -15852.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1703.8462
$ws.Range("I31").Value = 1452
$ws.Range("J31").Value = 8000
$ws.Range("K31").Value = 1452
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = -1157
$ws.Range("N31").Value = -8590
$ws.Range("H34").Value = 1703.8462
$ws.Range("I34").Value = 1452
$ws.Range("J34").Value = 8000
$ws.Range("K34").Value = 1452
$ws.Range("L34").Value = 8000
$ws.Range("M34").Value = -1250
$ws.Range("N34").Value = -8404
$ws.Range("H92").Value = 24841.834
$ws.Range("J92").Value = 24841.834
$ws.Range("L92").Value = 24841.834
$ws.Range("N92").Value = -29833.834
$ws.Range("H132").Value = 388061.5
$ws.Range("I132").Value = 501842.22
$ws.Range("J132").Value = 4051.5
$ws.Range("K132").Value = 1505526.66
$ws.Range("L132").Value = 12154.5
$ws.Range("M132").Value = -1502996.66
$ws.Range("N132").Value = -17214.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 2487
$ws.Range("I9").Value = 1999
$ws.Range("J9").Value = 2568.3333
$ws.Range("K9").Value = 5997
$ws.Range("L9").Value = 7704.999899999999
$ws.Range("M9").Value = -5773
$ws.Range("N9").Value = -8152.999899999999
$ws.Range("H70").Value = 11582.4
$ws.Range("I70").Value = 17304
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 51912
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -51597
$ws.Range("N70").Value = -9630
$ws.Range("H73").Value = 11582.4
$ws.Range("I73").Value = 17304
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 51912
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -50820
$ws.Range("N73").Value = -11184
$ws.Range("H80").Value = 6766.6665
$ws.Range("I80").Value = 16500
$ws.Range("J80").Value = 3985.7144
$ws.Range("K80").Value = 49500
$ws.Range("L80").Value = 11957.1432
$ws.Range("M80").Value = -48564
$ws.Range("N80").Value = -13829.1432
$ws.Range("H83").Value = 6766.6665
$ws.Range("I83").Value = 16500
$ws.Range("J83").Value = 3985.7144
$ws.Range("K83").Value = 148500
$ws.Range("L83").Value = 35871.4296
$ws.Range("M83").Value = -143820
$ws.Range("N83").Value = -45231.4296
$ws.Range("H87").Value = 8349
$ws.Range("I87").Value = 1675.375
$ws.Range("J87").Value = 9717.948
$ws.Range("K87").Value = 5026.125
$ws.Range("L87").Value = 29153.844
$ws.Range("M87").Value = -3778.125
$ws.Range("N87").Value = -31649.844
$ws.Range("H90").Value = 8349
$ws.Range("I90").Value = 1675.375
$ws.Range("J90").Value = 9717.948
$ws.Range("K90").Value = 15078.375
$ws.Range("L90").Value = 87461.53200000001
$ws.Range("M90").Value = -8838.375
$ws.Range("N90").Value = -99941.53200000001
$ws.Range("H122").Value = 1502
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1502
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13518
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -18418
$ws.Range("H131").Value = 25001378
$ws.Range("I131").Value = 439.0909
$ws.Range("J131").Value = 34484492
$ws.Range("K131").Value = 1317.2727
$ws.Range("L131").Value = 103453476
$ws.Range("M131").Value = 3722.7273
$ws.Range("N131").Value = -103463556
$ws.Range("H132").Value = 1970.9
$ws.Range("I132").Value = 1251
$ws.Range("J132").Value = 2450.8333
$ws.Range("K132").Value = 11259
$ws.Range("L132").Value = 22057.4997
$ws.Range("M132").Value = -8729
$ws.Range("N132").Value = -27117.4997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 325916.6
$ws.Range("I2").Value = 918307.0600000001
$ws.Range("J2").Value = 101.8
$ws.Range("K2").Value = 918307.0600000001
$ws.Range("L2").Value = 101.8
$ws.Range("M2").Value = -918194.0600000001
$ws.Range("N2").Value = -327.8
$ws.Range("H92").Value = 6343.4287
$ws.Range("J92").Value = 6343.4287
$ws.Range("L92").Value = 6343.4287
$ws.Range("N92").Value = -10087.4287
$ws.Range("H95").Value = 1443223.9
$ws.Range("J95").Value = 1443223.9
$ws.Range("L95").Value = 1443223.9
$ws.Range("N95").Value = -1448715.9
$ws.Range("H97").Value = 74661.57000000001
$ws.Range("I97").Value = 128657.75
$ws.Range("J97").Value = 2666.6667
$ws.Range("K97").Value = 128657.75
$ws.Range("L97").Value = 2666.6667
$ws.Range("M97").Value = -128161.75
$ws.Range("N97").Value = -3658.6667
$ws.Range("H113").Value = 5000
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 14666.667
$ws.Range("J94").Value = 14666.667
$ws.Range("L94").Value = 14666.667
$ws.Range("N94").Value = -16018.667
$ws.Range("H104").Value = 21173.076
$ws.Range("J104").Value = 21173.076
$ws.Range("L104").Value = 21173.076
$ws.Range("N104").Value = -28161.076
$ws.Range("H140").Value = 36142
$ws.Range("J140").Value = 36142
$ws.Range("L140").Value = 36142
$ws.Range("N140").Value = -46502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 9342.286
$ws.Range("J101").Value = 9342.286
$ws.Range("L101").Value = 9342.286
$ws.Range("N101").Value = -15832.286

Write-Host "Applied all Asura_Profits edits"